# fix(precios): corrige getPrice para retornar precio directo sin acceso a lista de precios
# - Elimina acceso innecesario a precios[priceList]
# - Retorna valor directo de product.precio
# - Los precios ya estan filtrados por lista al cargarlos
#
# Workbook-level consequence of the fix: the stray/duplicate reference value
# "evol115la" in GRUPO_SINPAR!A8 is cleared out, and the active sheet/tab
# selection moves from GRUPO_MAYORISTAS back to GRUPO_SINPAR.

$wb = $excel.ActiveWorkbook

$wsSinpar = $wb.Worksheets.Item("GRUPO_SINPAR")

# Clear the stale value out of GRUPO_SINPAR!A8 (removes the now-unused
# "evol115la" shared string as a side effect).
$wsSinpar.Range("A8").ClearContents()

# Make GRUPO_SINPAR the active sheet/tab again, with A8 selected.
# (GRUPO_MAYORISTAS keeps its prior C5 selection - left untouched here so it
# is no longer the active/tabSelected sheet.)
$wsSinpar.Activate()
$wsSinpar.Range("A8").Select()
